# Update excess mortality analyses and plots
#
# CBS weekly mortality data refresh:
#  - a handful of already-present weeks get revised "Waargenomen" (observed,
#    column G) counts (StatLine revisions), which ripple into the I column
#    (Oversterfte = G-H) through the existing shared formula;
#  - two new weeks (39 and 40) are appended to the table in rows 31-32;
#  - the "Som week 11 tot en met 19" totals row, previously directly under
#    the table at row 32, is pushed down to row 35 to leave a couple of
#    blank rows between the data and the totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Revised "Waargenomen" values for already-present weeks (col G) ---
# Column I (Oversterfte) recalculates automatically via the existing
# shared formula G{row}-H{row}.
$ws.Range("G11").Value = 2985
$ws.Range("G12").Value = 2774
$ws.Range("G17").Value = 2693
$ws.Range("G20").Value = 2617
$ws.Range("G22").Value = 2671
$ws.Range("G24").Value = 2635
$ws.Range("G25").Value = 3204
$ws.Range("G26").Value = 2842
$ws.Range("G27").Value = 2725
$ws.Range("G28").Value = 2672
$ws.Range("G29").Value = 2727
$ws.Range("G30").Value = 2694

# --- Move the totals row (F/G/H/I) from row 32 down to row 35, clearing ---
# --- the old row so rows 31-32 are free for the two new weeks below.    ---
$ws.Range("F35").Value = $ws.Range("F32").Value2
$ws.Range("G35").Formula = $ws.Range("G32").Formula
$ws.Range("H35").Formula = $ws.Range("H32").Formula
$ws.Range("I35").Formula = $ws.Range("I32").Formula
$ws.Range("G35:I35").NumberFormat = "0"

$ws.Range("F32:I32").Clear() | Out-Null

# --- New week 39 (row 31) ---
$ws.Range("F31").Value = 39
$ws.Range("G31").Value = 2865
$ws.Range("H31").Value = 2752
$ws.Range("I31").Formula = "=G31-H31"

# --- New week 40 (row 32) ---
$ws.Range("F32").Value = 40
$ws.Range("G32").Value = 2998
$ws.Range("H32").Value = 2786
$ws.Range("I32").Formula = "=G32-H32"

# Move the active selection, matching where the author last clicked.
$ws.Range("F33").Select()
